# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# for rows 2-20, reflecting the refreshed counts recorded in the commit.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2  = 9
    3  = 1086
    4  = 125
    5  = 86
    6  = 52
    7  = 57
    8  = 11212
    9  = 4286
    10 = 26
    11 = 24
    12 = 15
    13 = 2500
    14 = 1069
    15 = 102
    16 = 16
    17 = 160
    18 = 486
    19 = 11230
    20 = 11078
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
